$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 979.0833
$ws.Range("I8").Value = 38.77778
$ws.Range("J8").Value = 3800
$ws.Range("K8").Value = 116.33334
$ws.Range("L8").Value = 11400
$ws.Range("M8").Value = 22.66666000000001
$ws.Range("N8").Value = -11678

$ws.Range("H62").Value = 3624
$ws.Range("I62").Value = 2405
$ws.Range("J62").Value = 8500
$ws.Range("K62").Value = 2405
$ws.Range("L62").Value = 8500
$ws.Range("M62").Value = -1781
$ws.Range("N62").Value = -9748

$ws.Range("H65").Value = 3624
$ws.Range("I65").Value = 2405
$ws.Range("J65").Value = 8500
$ws.Range("K65").Value = 12025
$ws.Range("L65").Value = 42500
$ws.Range("M65").Value = -8905
$ws.Range("N65").Value = -48740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 549.9091
$ws.Range("I2").Value = 545.65625
$ws.Range("J2").Value = 561.25
$ws.Range("K2").Value = 545.65625
$ws.Range("L2").Value = 561.25
$ws.Range("M2").Value = -432.65625
$ws.Range("N2").Value = -787.25

$ws.Range("H32").Value = 7925.964
$ws.Range("I32").Value = 5842.5635
$ws.Range("K32").Value = 5842.5635
$ws.Range("M32").Value = -5555.5635

$ws.Range("H116").Value = 549.9091
$ws.Range("I116").Value = 545.65625
$ws.Range("J116").Value = 561.25
$ws.Range("K116").Value = 545.65625
$ws.Range("L116").Value = 561.25
$ws.Range("M116").Value = 1748.34375
$ws.Range("N116").Value = -5149.25

$ws.Range("H137").Value = 53780
$ws.Range("J137").Value = 53780
$ws.Range("L137").Value = 53780
$ws.Range("N137").Value = -63980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 549.9091
$ws.Range("I3").Value = 545.65625
$ws.Range("J3").Value = 561.25
$ws.Range("K3").Value = 545.65625
$ws.Range("L3").Value = 561.25
$ws.Range("M3").Value = -431.65625
$ws.Range("N3").Value = -789.25

$ws.Range("H107").Value = 508.7647
$ws.Range("I107").Value = 525.5
$ws.Range("K107").Value = 525.5
$ws.Range("M107").Value = 1394.5

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H134").Value = 2937.4792
$ws.Range("I134").Value = 1712.4872
$ws.Range("J134").Value = 8245.777
$ws.Range("K134").Value = 5137.461600000001
$ws.Range("L134").Value = 24737.331
$ws.Range("M134").Value = -2602.461600000001
$ws.Range("N134").Value = -29807.331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2033.5758
$ws.Range("I31").Value = 854.19446
$ws.Range("J31").Value = 3448.8333
$ws.Range("K31").Value = 854.19446
$ws.Range("L31").Value = 3448.8333
$ws.Range("M31").Value = -559.19446
$ws.Range("N31").Value = -4038.8333

$ws.Range("H34").Value = 2033.5758
$ws.Range("I34").Value = 854.19446
$ws.Range("J34").Value = 3448.8333
$ws.Range("K34").Value = 854.19446
$ws.Range("L34").Value = 3448.8333
$ws.Range("M34").Value = -652.19446
$ws.Range("N34").Value = -3852.8333

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 534.75
$ws.Range("I113").Value = 549.61365
$ws.Range("J113").Value = 507.5
$ws.Range("K113").Value = 1648.84095
$ws.Range("L113").Value = 1522.5
$ws.Range("M113").Value = 521.15905
$ws.Range("N113").Value = -5862.5

$ws.Range("H122").Value = 2329.2825
$ws.Range("J122").Value = 3095.4834
$ws.Range("L122").Value = 27859.3506
$ws.Range("N122").Value = -32759.3506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2733.639
$ws.Range("I102").Value = 2205.4285
$ws.Range("J102").Value = 3473.1333
$ws.Range("K102").Value = 2205.4285
$ws.Range("L102").Value = 3473.1333
$ws.Range("M102").Value = -583.4285
$ws.Range("N102").Value = -6717.1333

$ws.Range("H107").Value = 494.3846
$ws.Range("I107").Value = 339
$ws.Range("J107").Value = 844
$ws.Range("K107").Value = 339
$ws.Range("L107").Value = 844
$ws.Range("M107").Value = 1581
$ws.Range("N107").Value = -4684

$ws.Range("H113").Value = 1408.1666
$ws.Range("I113").Value = 1389.8
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1389.8
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 780.2
$ws.Range("N113").Value = -5840

$ws.Range("H132").Value = 1744.6538
$ws.Range("I132").Value = 1043.878
$ws.Range("J132").Value = 4356.636
$ws.Range("K132").Value = 3131.634
$ws.Range("L132").Value = 13069.908
$ws.Range("M132").Value = -601.634
$ws.Range("N132").Value = -18129.908

$ws.Range("H136").Value = 11336.667
$ws.Range("J136").Value = 11336.667
$ws.Range("L136").Value = 34010.001
$ws.Range("N136").Value = -39110.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4451.56
$ws.Range("I40").Value = 4362.1577
$ws.Range("J40").Value = 4734.6665
$ws.Range("K40").Value = 4362.1577
$ws.Range("L40").Value = 4734.6665
$ws.Range("M40").Value = -4226.1577
$ws.Range("N40").Value = -5006.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 35796200
$ws.Range("I62").Value = 71432136
$ws.Range("J62").Value = 160257.58
$ws.Range("K62").Value = 71432136
$ws.Range("L62").Value = 160257.58
$ws.Range("M62").Value = -71431512
$ws.Range("N62").Value = -161505.58

$ws.Range("H65").Value = 35796200
$ws.Range("I65").Value = 71432136
$ws.Range("J65").Value = 160257.58
$ws.Range("K65").Value = 357160680
$ws.Range("L65").Value = 801287.8999999999
$ws.Range("M65").Value = -357157560
$ws.Range("N65").Value = -807527.8999999999

$ws.Range("H122").Value = 2421.2083
$ws.Range("I122").Value = 1755.4062
$ws.Range("J122").Value = 3752.8125
$ws.Range("K122").Value = 5266.2186
$ws.Range("L122").Value = 11258.4375
$ws.Range("M122").Value = -2816.2186
$ws.Range("N122").Value = -16158.4375

$ws.Range("H126").Value = 416474.3
$ws.Range("I126").Value = 2700.5334
$ws.Range("K126").Value = 8101.600199999999
$ws.Range("M126").Value = -5631.600199999999

$ws.Range("H132").Value = 5953920.5
$ws.Range("I132").Value = 532.4643
$ws.Range("J132").Value = 11907308
$ws.Range("K132").Value = 1597.3929
$ws.Range("L132").Value = 35721924
$ws.Range("M132").Value = 932.6071000000002
$ws.Range("N132").Value = -35726984
